$d = $word.ActiveDocument

# 1) "dorostenci " + "(15-" -> single run "dorostenci (15-"
#    (also coalesces the following "19 let" + ")" runs into "19 let)")
$d.Content.Find.Execute(
    "dorostenci (15-", $true, $false, $false, $false, $false,
    $true, 1, $false, "dorostenci (15-", 2) | Out-Null

# 2) Update the closing motivational sentence.
$d.Content.Find.Execute(
    "Není důležité zvítězit, ale zúčastnit se a zůstat zdravý!", $true, $false, $false, $false, $false,
    $true, 1, $false, "Není důležité zvítězit, ale zůstat v kondici a pohodě.", 2) | Out-Null
